{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfunction setCellText(rowIndex, text) {\n  const cell = table.getCell(rowIndex, 0);\n  const range = cell.body.getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\n// 1) First three summary rows: 99.99 / 0.05 / 834 -> 0M\nsetCellText(0, \"0M\");\nsetCellText(1, \"0M\");\nsetCellText(2, \"0M\");\n\n// 2) Fourth row: 3 -> 1042\nsetCellText(3, \"1042\");\nawait context.sync();\n\n// 3) Insert two new rows right after row 4 (index 3), carrying the new\n//    per-iteration values; they pick up the surrounding cell formatting\n//    (Times New Roman, sz 22) automatically.\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst anchorRow = rows.items[4]; // currently \"0.00004\"\nanchorRow.insertRows(Word.InsertLocation.before, 2, [[\"0.00002\"], [\"0.00018\"]]);\nawait context.sync();\n\n// Row \"0.00004\" (now index 6) is unchanged; the next two rows (both\n// \"0.00005\", now indices 7 and 8) are removed.\nrows.load(\"items\");\nawait context.sync();\nrows.items[7].delete();\nrows.items[7].delete();\nawait context.sync();\n\n// 4) Row with 0.00014 -> 0.04505 (still at index 11 since two rows were\n//    added and two removed ahead of it)\nsetCellText(11, \"0.04505\");\nawait context.sync();\n\n// 5) Collapse the final three multi-column summary rows down to their\n//    first value, dropping the tab-separated remainder.\nsetCellText(43, \"99.99\");\nsetCellText(44, \"0.05\");\nsetCellText(45, \"834\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1) First three summary rows: 99.99 / 0.05 / 834 -> 0M\n$t.Rows.Item(1).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(2).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(3).Cells.Item(1).Range.Text = \"0M\"\n\n# 2) Fourth row: 3 -> 1042\n$t.Rows.Item(4).Cells.Item(1).Range.Text = \"1042\"\n\n# 3) Insert two new rows before the \"0.00004\" row (currently row 5), carrying\n#    the new per-iteration values; Rows.Add(beforeRow) inserts ahead of the target\n$newRow1 = $t.Rows.Add($t.Rows.Item(5))\n$newRow1.Cells.Item(1).Range.Text = \"0.00002\"\n\n$newRow2 = $t.Rows.Add($t.Rows.Item(6))\n$newRow2.Cells.Item(1).Range.Text = \"0.00018\"\n\n# Row 7 (\"0.00004\") is unchanged; the next two rows (both \"0.00005\") are removed\n$t.Rows.Item(8).Delete()\n$t.Rows.Item(8).Delete()\n\n# 4) Row with 0.00014 -> 0.04505\n$t.Rows.Item(12).Cells.Item(1).Range.Text = \"0.04505\"\n\n# 5) Collapse the final three multi-column summary rows down to their first value\n$t.Rows.Item(44).Cells.Item(1).Range.Text = \"99.99\"\n$t.Rows.Item(45).Cells.Item(1).Range.Text = \"0.05\"\n$t.Rows.Item(46).Cells.Item(1).Range.Text = \"834\"\n"}
